# 4_DefinitionenDerProzesswoerter.xlsx - update process-word table:
#  - remove the "modifizieren" entry (row whose word is "modifizieren")
#  - append 7 newly defined process words (no definition text yet) at the
#    bottom of the table
#  - refresh the view/selection state to match the authored workbook

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- locate & delete the "modifizieren" row -------------------------------
$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$targetRow = 0
for ($r = 1; $r -le $lastRow; $r++) {
    $word = $ws.Cells.Item($r, 1).Value2
    if ($word -eq "modifizieren") {
        $targetRow = $r
        break
    }
}
if ($targetRow -gt 0) {
    $ws.Rows($targetRow).Delete()
}

# --- append the newly introduced process words (column A only) ----------
$newWords = @("aktivieren", "deaktivieren", "priorisieren", "behandeln", "ignorieren", "pausieren", "fortsetzen")

$lastUsedRow = $ws.UsedRange.Rows.Count
$nextRow = $lastUsedRow + 1
foreach ($word in $newWords) {
    $ws.Cells.Item($nextRow, 1).Value = $word
    $nextRow++
}

# --- refresh sort bookkeeping over the (now 17-word) definition table ----
$lastDataRow = $nextRow - [int]$newWords.Length - 1
$sortObj = $ws.Sort
$sortObj.SetRange($ws.Range("A2:B$lastDataRow"))
$sortObj.Header = 2
$sortObj.Apply()

# --- restore view state (selection / scroll position) --------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 8
$win.ScrollColumn = 1
$ws.Range("D1").Select()
